$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "302.72"
Set-TextValue "E2" "4.56%"
Set-TextValue "D3" "35.79"
Set-TextValue "E3" "15.21%"
Set-TextValue "D4" "5.173"
Set-TextValue "E4" "4.48%"
Set-TextValue "D5" "0.07849"
Set-TextValue "E5" "6.70%"
Set-TextValue "D6" "2.284"
Set-TextValue "E6" "-1.06%"
Set-TextValue "D7" "8.058"
Set-TextValue "E7" "4.94%"
Set-TextValue "D8" "3.970"
Set-TextValue "E8" "6.03%"
Set-TextValue "D9" "0.9259"
Set-TextValue "E9" "0.83%"
Set-TextValue "D10" "0.1013"
Set-TextValue "E10" "10.76%"
Set-TextValue "D11" "0.1831"
Set-TextValue "E11" "7.32%"
Set-TextValue "D12" "0.08528"
Set-TextValue "E12" "4.59%"
Set-TextValue "D13" "0.03388"
Set-TextValue "E13" "8.84%"
Set-TextValue "D14" "0.09923"
Set-TextValue "E14" "-0.72%"
Set-TextValue "D15" "0.001473"
Set-TextValue "E15" "-1.73%"
Set-TextValue "D16" "0.005750"
Set-TextValue "E16" "0.29%"
Set-TextValue "D17" "3.483"
Set-TextValue "E17" "0.32%"
Set-TextValue "E18" "4.79%"
Set-TextValue "D19" "0.3431"
Set-TextValue "E19" "3.00%"
Set-TextValue "D20" "0.1324"
Set-TextValue "E20" "1.93%"
Set-TextValue "D21" "4.532"
Set-TextValue "E21" "8.40%"
Set-TextValue "D22" "0.2214"
Set-TextValue "E22" "4.25%"
Set-TextValue "D23" "0.04625"
Set-TextValue "E23" "2.52%"
Set-TextValue "D24" "0.001215"
Set-TextValue "E24" "0.00%"
Set-TextValue "D25" "0.004491"
Set-TextValue "E25" "6.96%"
Set-TextValue "D26" "0.0001294"
Set-TextValue "E26" "-0.48%"
Set-TextValue "D27" "0.0003383"
Set-TextValue "E27" "-0.32%"
Set-TextValue "D39" "0.01746"
Set-TextValue "E39" "10.55%"
Set-TextValue "D40" "0.04729"
Set-TextValue "E40" "4.84%"
Set-TextValue "D41" "0.007835"
Set-TextValue "E41" "5.98%"
Set-TextValue "D42" "0.1415"
Set-TextValue "E42" "5.83%"
Set-TextValue "D43" "0.008779"
Set-TextValue "E43" "-10.88%"
Set-TextValue "D44" "0.002208"
Set-TextValue "E44" "-0.55%"
Set-TextValue "D45" "0.009147"
Set-TextValue "E45" "7.55%"
Set-TextValue "D46" "0.00006032"
Set-TextValue "E46" "-1.19%"
Set-TextValue "E47" "-0.43%"
Set-TextValue "D48" "5.804"
Set-TextValue "E48" "126.27%"
Set-TextValue "D49" "0.002678"
Set-TextValue "E49" "33.91%"
Set-TextValue "D50" "0.00002091"
Set-TextValue "E50" "-0.43%"
Set-TextValue "D51" "0.0001991"
Set-TextValue "E51" "-0.43%"
